$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear example/placeholder note lines under "Catatan Hutang 1" so the
# exported template starts blank instead of carrying sample data.
$ws.Range("B22").Value = ""
$ws.Range("B23").Value = ""
$ws.Range("B24").Value = ""

# Clear example/placeholder note lines under "Catatan Pembahasan 2".
$ws.Range("B28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("B30").Value = ""
$ws.Range("B31").Value = ""
$ws.Range("B32").Value = ""
$ws.Range("B33").Value = ""
$ws.Range("B34").Value = ""

# Clear example/placeholder note lines + amounts under "Catatan Beban HO 3".
$ws.Range("B38").Value = ""
$ws.Range("C38").Value = ""
$ws.Range("B39").Value = ""
$ws.Range("C39").Value = ""

# Clear the sample total amount.
$ws.Range("C40").Value = ""

# Leave the cursor on G25 and scroll the sheet back to the top, matching
# the template's saved view state.
$ws.Range("G25").Select()
